$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1:X1").ClearContents()
